$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 15 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 15
}

# Update the Actual Production (MW) values for rows 28-31 in column B
$ws.Cells.Item(28, 2).Value2 = 9
$ws.Cells.Item(29, 2).Value2 = 34
$ws.Cells.Item(30, 2).Value2 = 70
$ws.Cells.Item(31, 2).Value2 = 112
